$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the woe values for registered_via rows 2-6
$ws.Range("B2").Value = -1.052
$ws.Range("B3").Value = 0.329
$ws.Range("B4").Value = 0.625
$ws.Range("B5").Value = 0.6840000000000001
$ws.Range("B6").Value = 1.125

# Delete row 7 entirely (A7=10, B7 empty) so the used range shrinks to A1:B6
$ws.Rows.Item(7).Delete()
